# ---------------------------------------------------------------------------
# Edit summary (per the target diff):
#   1. Every table's <a:tableStyleId> is changed from
#      {6C14C857-4E22-4F60-A55F-FC2D71E598A1} to
#      {BA355E6F-4B2B-45DA-94FA-C7ED8F1C87BC}.
#   2. The theme color palette that the deck actually renders with (the one
#      shared by the slide master / presentation, persisted as
#      ppt/theme/theme2.xml) is swapped from the custom "NYT Slides
#      Template" ("Simple Light") palette to the plain "Default" Office
#      palette that used to live in ppt/theme/theme1.xml.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck -----------------------------------
$newTableStyleId = "{BA355E6F-4B2B-45DA-94FA-C7ED8F1C87BC}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2. Swap the rendered theme's color scheme back to the plain "Default"
#        Office palette (was previously the "NYT Slides Template" palette).
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # dk1
$tcs.Item(2).RGB  = 0xFFFFFF   # lt1
$tcs.Item(3).RGB  = 0x588115   # dk2      (158158 stored BGR)
$tcs.Item(4).RGB  = 0xF3F3F3   # lt2
$tcs.Item(5).RGB  = 0xC78D05   # accent1  (058DC7 stored BGR)
$tcs.Item(6).RGB  = 0x32B450   # accent2  (50B432 stored BGR)
$tcs.Item(7).RGB  = 0x1B56ED   # accent3  (ED561B stored BGR)
$tcs.Item(8).RGB  = 0x00EFED   # accent4  (EDEF00 stored BGR)
$tcs.Item(9).RGB  = 0xE5CB24   # accent5  (24CBE5 stored BGR)
$tcs.Item(10).RGB = 0x72E564   # accent6  (64E572 stored BGR)
$tcs.Item(11).RGB = 0xCC0022   # hlink    (2200CC stored BGR)
$tcs.Item(12).RGB = 0x8B1A55   # folHlink (551A8B stored BGR)
